$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Remove the duplicate "Contact" row (old row 11 duplicated old row 10).
# This shifts every following row up by one (old row 12 -> new row 11, etc.)
$ws.Rows.Item(11).Delete()

# Version bump
$ws.Range("B3").Value = "6.0.0"

# Updated publication date
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was previously blank)
$ws.Range("B9").Value = "Alvearie Team"

# Replace the leftover "Contact" row with a new "Jurisdiction" property/value pair
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# "Case Sensitive" row gets a value of text "true" (not boolean TRUE).
# Excel's Range.Value setter auto-coerces the literal "true"/"false" text into
# a native Boolean, so build the text value via a formula and paste it back in
# as a value (Copy + PasteSpecial values-only) to keep it as a real text/string cell.
$ws.Range("B12").Formula = "=""true"""
$ws.Range("B12").Copy()
$ws.Range("B14").PasteSpecial(-4163, 0, $false, $false)
$excel.CutCopyMode = 0
$ws.Range("B12").ClearContents()
